# The workbook tracks daily price observations for "Arveja Verde" at the
# Vega Central Mapocho de Santiago market. A new daily record was inserted
# as row 48, pushing all subsequent rows (old 48..129) down by one
# (new 49..130).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 48, shifting existing rows 48-129 down to 49-130.
$ws.Rows.Item(48).Insert()

# Populate the newly inserted row 48 with the new observation.
$ws.Cells.Item(48, 1).Value = 9
$ws.Cells.Item(48, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(48, 3).Value = "Metropolitana"
$ws.Cells.Item(48, 4).Value = 44797
$ws.Cells.Item(48, 5).Value = 13
$ws.Cells.Item(48, 6).Value = 100112022
$ws.Cells.Item(48, 7).Value = "Arveja Verde"
$ws.Cells.Item(48, 8).Value = "Sin especificar"
$ws.Cells.Item(48, 9).Value = "Primera"
$ws.Cells.Item(48, 10).Value = 34
$ws.Cells.Item(48, 11).Value = 40000
$ws.Cells.Item(48, 12).Value = 40000
$ws.Cells.Item(48, 13).Value = 40000
$ws.Cells.Item(48, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(48, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(48, 16).Value = 1600
$ws.Cells.Item(48, 17).Value = 25
$ws.Cells.Item(48, 18).Value = "Hortaliza"
